# AFDP-2522 ensure case files rules work before queue has been assigned
#
# The "Save Case File Rules" decision table (Sheet1) uses Spring-EL
# condition strings like `queue.name == 'Billing'` to decide whether to
# stamp/clear the Billing/Hold "enter date" fields. When a case file has
# not yet been assigned to a queue, `queue` is null and the plain `.name`
# navigation throws instead of short-circuiting. Switch those four
# conditions to the null-safe navigation operator (`queue?.name`) so the
# rules evaluate safely even before a queue has been assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C28").Value = "queue?.name == 'Billing' && billingEnterDate == null"
$ws.Range("C29").Value = "queue?.name != 'Billing'"
$ws.Range("C30").Value = "queue?.name == 'Hold' && holdEnterDate == null"
$ws.Range("C31").Value = "queue?.name != 'Hold'"

# Trim the trailing blank rows left over below the table.
$ws.Range("B32:D32").Clear()
$ws.Rows.Item(33).Delete()

# Leave the selection where the author left it after the edit.
[void]$ws.Range("B26").Select()
